$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 3006362.2
$ws.Cells.Item(15, 9).Value = 3006362.2
$ws.Cells.Item(15, 11).Value = 9019086.600000001
$ws.Cells.Item(15, 13).Value = -9018917.600000001
$ws.Cells.Item(17, 8).Value = 1812.5
$ws.Cells.Item(17, 10).Value = 1812.5
$ws.Cells.Item(17, 12).Value = 5437.5
$ws.Cells.Item(17, 14).Value = -5773.5
$ws.Cells.Item(32, 8).Value = 3493.0625
$ws.Cells.Item(32, 10).Value = 3367.9
$ws.Cells.Item(32, 12).Value = 3367.9
$ws.Cells.Item(32, 14).Value = -4019.9
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).ClearContents()
$ws.Cells.Item(76, 8).Value = 5665
$ws.Cells.Item(76, 9).Value = 4258.5713
$ws.Cells.Item(76, 11).Value = 4258.5713
$ws.Cells.Item(76, 13).Value = -3943.5713
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).ClearContents()
$ws.Cells.Item(79, 8).Value = 5665
$ws.Cells.Item(79, 9).Value = 4258.5713
$ws.Cells.Item(79, 11).Value = 4258.5713
$ws.Cells.Item(79, 13).Value = -3166.5713
$ws.Cells.Item(112, 8).Value = 1739.7
$ws.Cells.Item(112, 10).Value = 1912.25
$ws.Cells.Item(112, 12).Value = 5736.75
$ws.Cells.Item(112, 14).Value = -7952.75
$ws.Cells.Item(129, 8).Value = 1366.3334
$ws.Cells.Item(129, 9).Value = 818.7273
$ws.Cells.Item(129, 10).Value = 2872.25
$ws.Cells.Item(129, 11).Value = 2456.1819
$ws.Cells.Item(129, 12).Value = 8616.75
$ws.Cells.Item(129, 13).Value = 2543.8181
$ws.Cells.Item(129, 14).Value = -18616.75
$ws.Cells.Item(138, 8).Value = 1856.3
$ws.Cells.Item(138, 10).Value = 2339.7285
$ws.Cells.Item(138, 12).Value = 7019.185500000001
$ws.Cells.Item(138, 14).Value = -17299.1855
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 354.92307
$ws.Cells.Item(5, 9).Value = 197.5
$ws.Cells.Item(5, 10).Value = 606.8
$ws.Cells.Item(5, 11).Value = 197.5
$ws.Cells.Item(5, 12).Value = 606.8
$ws.Cells.Item(5, 13).Value = -85.5
$ws.Cells.Item(5, 14).Value = -830.8
$ws.Cells.Item(32, 8).Value = 10874884
$ws.Cells.Item(32, 9).Value = 12822724
$ws.Cells.Item(32, 10).Value = 22635.285
$ws.Cells.Item(32, 11).Value = 12822724
$ws.Cells.Item(32, 12).Value = 22635.285
$ws.Cells.Item(32, 13).Value = -12822437
$ws.Cells.Item(32, 14).Value = -23209.285
$ws.Cells.Item(110, 8).Value = 2076.7715
$ws.Cells.Item(110, 9).Value = 1863
$ws.Cells.Item(110, 11).Value = 1863
$ws.Cells.Item(110, 13).Value = 182
$ws.Cells.Item(120, 8).Value = 73240
$ws.Cells.Item(120, 10).Value = 73240
$ws.Cells.Item(120, 12).Value = 73240
$ws.Cells.Item(120, 14).Value = -82916
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 354.92307
$ws.Cells.Item(4, 9).Value = 197.5
$ws.Cells.Item(4, 10).Value = 606.8
$ws.Cells.Item(4, 11).Value = 197.5
$ws.Cells.Item(4, 12).Value = 606.8
$ws.Cells.Item(4, 13).Value = -82.5
$ws.Cells.Item(4, 14).Value = -836.8
$ws.Cells.Item(100, 8).Value = 37000
$ws.Cells.Item(100, 10).Value = 37000
$ws.Cells.Item(100, 12).Value = 37000
$ws.Cells.Item(100, 14).Value = -39164
$ws.Cells.Item(112, 8).Value = 107994.5
$ws.Cells.Item(112, 10).Value = 107994.5
$ws.Cells.Item(112, 12).Value = 107994.5
$ws.Cells.Item(112, 14).Value = -110948.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 8780.799999999999
$ws.Cells.Item(89, 8).Value = 8780.799999999999
$ws.Cells.Item(134, 8).Value = 1959.78
$ws.Cells.Item(134, 9).Value = 1520.3658
$ws.Cells.Item(134, 10).Value = 3961.5557
$ws.Cells.Item(134, 11).Value = 4561.097400000001
$ws.Cells.Item(134, 12).Value = 11884.6671
$ws.Cells.Item(134, 13).Value = -2026.097400000001
$ws.Cells.Item(134, 14).Value = -16954.6671
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 237.5
$ws.Cells.Item(2, 9).Value = 459.2857
$ws.Cells.Item(2, 10).Value = 163.57143
$ws.Cells.Item(2, 11).Value = 2755.7142
$ws.Cells.Item(2, 12).Value = 981.42858
$ws.Cells.Item(2, 13).Value = -2642.7142
$ws.Cells.Item(2, 14).Value = -1207.42858
$ws.Cells.Item(5, 8).Value = 2024.1538
$ws.Cells.Item(5, 9).Value = 2035.75
$ws.Cells.Item(5, 10).Value = 1885
$ws.Cells.Item(5, 11).Value = 6107.25
$ws.Cells.Item(5, 12).Value = 5655
$ws.Cells.Item(5, 13).Value = -5995.25
$ws.Cells.Item(5, 14).Value = -5879
$ws.Cells.Item(109, 8).Value = 2619
$ws.Cells.Item(109, 9).Value = 2619
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 11).Value = 7857
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 13).Value = -6817
$ws.Cells.Item(109, 14).ClearContents()
$ws.Cells.Item(131, 8).Value = 5275.6
$ws.Cells.Item(131, 10).Value = 5275.6
$ws.Cells.Item(131, 12).Value = 15826.8
$ws.Cells.Item(131, 14).Value = -25906.8
$ws.Cells.Item(132, 8).Value = 1596.4
$ws.Cells.Item(132, 10).Value = 1745.5
$ws.Cells.Item(132, 12).Value = 15709.5
$ws.Cells.Item(132, 14).Value = -20769.5
$ws.Cells.Item(135, 8).Value = 2024.1538
$ws.Cells.Item(135, 9).Value = 2035.75
$ws.Cells.Item(135, 10).Value = 1885
$ws.Cells.Item(135, 11).Value = 18321.75
$ws.Cells.Item(135, 12).Value = 16965
$ws.Cells.Item(135, 13).Value = -15786.75
$ws.Cells.Item(135, 14).Value = -22035
$ws.Cells.Item(137, 8).Value = 8161
$ws.Cells.Item(137, 9).Value = 10015
$ws.Cells.Item(137, 10).Value = 7419.4
$ws.Cells.Item(137, 11).Value = 30045
$ws.Cells.Item(137, 12).Value = 22258.2
$ws.Cells.Item(137, 13).Value = -24945
$ws.Cells.Item(137, 14).Value = -32458.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(47, 8).Value = 14995
$ws.Cells.Item(47, 10).Value = 14995
$ws.Cells.Item(47, 12).Value = 14995
$ws.Cells.Item(47, 14).Value = -16131
$ws.Cells.Item(126, 8).Value = 3866.2778
$ws.Cells.Item(126, 9).Value = 3428.4443
$ws.Cells.Item(126, 11).Value = 10285.3329
$ws.Cells.Item(126, 13).Value = -7815.332900000001
$ws.Cells.Item(132, 8).Value = 20413784
$ws.Cells.Item(132, 9).Value = 35717340
$ws.Cells.Item(132, 10).Value = 9043
$ws.Cells.Item(132, 11).Value = 107152020
$ws.Cells.Item(132, 12).Value = 27129
$ws.Cells.Item(132, 13).Value = -107149490
$ws.Cells.Item(132, 14).Value = -32189
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1336.5454
$ws.Cells.Item(22, 9).Value = 1300
$ws.Cells.Item(22, 10).Value = 1380.4
$ws.Cells.Item(22, 11).Value = 1300
$ws.Cells.Item(22, 12).Value = 1380.4
$ws.Cells.Item(22, 13).Value = -1005
$ws.Cells.Item(22, 14).Value = -1970.4
$ws.Cells.Item(27, 8).Value = 1336.5454
$ws.Cells.Item(27, 9).Value = 1300
$ws.Cells.Item(27, 10).Value = 1380.4
$ws.Cells.Item(27, 11).Value = 1300
$ws.Cells.Item(27, 12).Value = 1380.4
$ws.Cells.Item(27, 13).Value = -1193
$ws.Cells.Item(27, 14).Value = -1594.4
$ws.Cells.Item(55, 8).Value = 47619748
$ws.Cells.Item(55, 10).Value = 668.1667
$ws.Cells.Item(55, 12).Value = 668.1667
$ws.Cells.Item(55, 14).Value = -1014.1667
$ws.Cells.Item(61, 8).Value = 1386.9231
$ws.Cells.Item(61, 9).Value = 1247.6
$ws.Cells.Item(61, 10).Value = 1851.3334
$ws.Cells.Item(61, 11).Value = 1247.6
$ws.Cells.Item(61, 12).Value = 1851.3334
$ws.Cells.Item(61, 13).Value = -1045.6
$ws.Cells.Item(61, 14).Value = -2255.3334
$ws.Cells.Item(96, 8).Value = 97500
$ws.Cells.Item(96, 10).Value = 97500
$ws.Cells.Item(96, 12).Value = 97500
$ws.Cells.Item(96, 14).Value = -102992
$ws.Cells.Item(109, 8).Value = 100074.664
$ws.Cells.Item(109, 10).Value = 100074.664
$ws.Cells.Item(109, 12).Value = 100074.664
$ws.Cells.Item(109, 14).Value = -102848.664
$ws.Cells.Item(113, 8).Value = 1386.9231
$ws.Cells.Item(113, 9).Value = 1247.6
$ws.Cells.Item(113, 10).Value = 1851.3334
$ws.Cells.Item(113, 11).Value = 1247.6
$ws.Cells.Item(113, 12).Value = 1851.3334
$ws.Cells.Item(113, 13).Value = 922.4000000000001
$ws.Cells.Item(113, 14).Value = -6191.3334
$ws.Cells.Item(122, 8).Value = 5711.516
$ws.Cells.Item(122, 9).Value = 5175.5713
$ws.Cells.Item(122, 11).Value = 15526.7139
$ws.Cells.Item(122, 13).Value = -13076.7139
$ws.Cells.Item(123, 8).Value = 55995
$ws.Cells.Item(123, 10).Value = 55995
$ws.Cells.Item(123, 12).Value = 55995
$ws.Cells.Item(123, 14).Value = -65795
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 3388.7778
$ws.Cells.Item(107, 8).Value = 11628735
$ws.Cells.Item(107, 9).Value = 15625920
$ws.Cells.Item(107, 11).Value = 46877760
$ws.Cells.Item(107, 13).Value = -46875840
$ws.Cells.Item(109, 8).Value = 107930
$ws.Cells.Item(109, 10).Value = 107930
$ws.Cells.Item(109, 12).Value = 107930
$ws.Cells.Item(109, 14).Value = -110704
$ws.Cells.Item(125, 8).Value = 42843.89
$ws.Cells.Item(125, 10).Value = 42843.89
$ws.Cells.Item(125, 12).Value = 42843.89
$ws.Cells.Item(125, 14).Value = -52683.89
